$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking price strings (e.g. "237.50", "1.000")
# are not auto-converted to numbers by Excel when assigned via .Value
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.251.28"
$ws.Range("D3").Value = "1.863.91"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "237.50"
$ws.Range("E5").Value = "  +1.61%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "0.4684"
$ws.Range("D8").Value = "0.2866"
$ws.Range("D9").Value = "0.06549"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "22.29"
$ws.Range("E10").Value = "  +12.25%  "
$ws.Range("D11").Value = "0.07894"
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("D12").Value = "97.91"
$ws.Range("E12").Value = "  +1.52%  "
$ws.Range("D13").Value = "1.867.10"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").Value = "5.190"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").Value = "0.6813"
$ws.Range("E15").Value = "  +2.23%  "
$ws.Range("D16").Value = "278.54"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "30.250.30"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "13.60"
$ws.Range("E18").Value = "  +7.80%  "
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "0.000007342"
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("D21").Value = "5.385"
$ws.Range("E21").Value = "  -2.41%  "
$ws.Range("D22").Value = "2.110.52"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "6.200"
$ws.Range("E24").Value = "  +1.20%  "
$ws.Range("D25").Value = "168.38"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("D26").Value = "9.278"
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("D27").Value = "19.09"
$ws.Range("E27").Value = "  +1.08%  "
$ws.Range("D28").Value = "1.947"
$ws.Range("E28").Value = "  +1.97%  "
$ws.Range("D29").Value = "1.382"
$ws.Range("E29").Value = "  +3.16%  "
$ws.Range("D30").Value = "0.09815"
$ws.Range("E30").Value = "  +2.52%  "
$ws.Range("D31").Value = "4.385"
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("E32").Value = "  +0.81%  "
$ws.Range("D33").Value = "4.074"
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("D34").Value = "0.04754"
$ws.Range("E34").Value = "  +2.09%  "
$ws.Range("D35").Value = "1.141"
$ws.Range("E35").Value = "  +4.52%  "
$ws.Range("D36").Value = "0.7049"
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("D37").Value = "2.704"
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").Value = "0.01880"
$ws.Range("E38").Value = "  +1.58%  "
$ws.Range("D39").Value = "2.623"
$ws.Range("E39").Value = "  +4.14%  "
$ws.Range("D40").Value = "76.35"
$ws.Range("E40").Value = "  +4.04%  "
$ws.Range("D41").Value = "6.290"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("D42").Value = "1.958"
$ws.Range("E42").Value = "  +2.07%  "
$ws.Range("D43").Value = "0.8511"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D44").Value = "0.4187"
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("D45").Value = "0.9998"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "103.15"
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("D47").Value = "7.226"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("D48").Value = "951.27"
$ws.Range("E48").Value = "  -3.95%  "
$ws.Range("D49").Value = "9.283"
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").Value = "34.28"
$ws.Range("D51").Value = "0.05638"
$ws.Range("E51").Value = "  -0.14%  "

# Restore default (unstyled) formatting on column D now that values are locked in as text
$ws.Range("D2:D51").Style = "Normal"

Write-Host "Updated cryptos list"
